# Applies the MUSE "data standardization" document edit:
#  - Removes the Heading3 style from the existing Q&A paragraphs (they become
#    plain body paragraphs).
#  - Inserts "for fuels " into the import/export price question.
#  - Rewrites the "Can MUSE deal with HFO and LFO..." question as a new
#    question about multiple Process outputs in the Technodata file.
#  - Changes "Do you have access" to "Do we have access".
#  - Appends a new "Data issues" (Heading1) section with a paragraph about
#    Eritrea's Off-grid solar PV data gap.

$d = $word.ActiveDocument

# --- 1. Strip the Heading3 style from the existing Q&A paragraphs -----------
# (paragraph indices are stable; empty separator paragraphs stay "Normal")
$qaIndexes = @(3, 5, 7, 9, 11, 13, 15)
foreach ($idx in $qaIndexes) {
    $d.Paragraphs($idx).Style = "Normal"
}

# --- 2. "How to deal with import/export price values ..." -> add "for fuels "
$d.Content.Find.Execute(
    "How to deal with import/export price values (I have import/export prices in `$/GJ)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "How to deal with import/export price values (I have import/export prices for fuels in `$/GJ)",
    2)

# --- 3. Replace the HFO/LFO question with the Technodata multiple-outputs Q.
$d.Content.Find.Execute(
    "Can MUSE deal with HFO and LFO output ratios from crude oil?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "How can MUSE deal with multiple outputs of a Process in the Technodata file? In this work Crude Oil Refineries produce HFO and LFO.",
    2)

# --- 4. "Do you have access" -> "Do we have access" -------------------------
$d.Content.Find.Execute(
    "Do you have access to typical",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Do we have access to typical",
    2)

# --- 5. Append the new "Data issues" section --------------------------------
# Target tail layout:  ... <empty> <Heading1 "Data issues"> <empty> <Normal Eritrea-paragraph>
$countBefore = $d.Paragraphs.Count

# empty separator paragraph (Normal)
$r = $d.Paragraphs($countBefore).Range
$r.Collapse(0)   # wdCollapseEnd
$r.InsertParagraphAfter()
$d.Paragraphs($countBefore + 1).Style = "Normal"

# "Data issues" heading paragraph
$r = $d.Paragraphs($countBefore + 1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$headingPara = $d.Paragraphs($countBefore + 2)
$headingPara.Range.Text = "Data issues"
$headingPara.Style = "Heading1"

# empty separator paragraph (Normal)
$r = $d.Paragraphs($countBefore + 2).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs($countBefore + 3).Style = "Normal"

# final paragraph about Eritrea / Off-grid solar PV data gap
$r = $d.Paragraphs($countBefore + 3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$finalPara = $d.Paragraphs($countBefore + 4)
$finalPara.Range.Text = "For Eritrea, for example, there is Off-grid solar PV in the raw existing capacity table (Table1), but no data for Off-grid solar PV in the raw Technodata table (Table2)."
$finalPara.Style = "Normal"

Write-Output "done"
